$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add H1 (04-09_A) and I1 (04-09_0), copying formats from F1/G1
$ws.Range("F1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "04-09_A"
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "04-09_0"

# Row 2
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4163)
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = 0
$ws.Range("G2").Value = 0

# Row 3
$ws.Range("G3").Copy()
$ws.Range("I3").PasteSpecial(-4163)
$ws.Range("F3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = 0
$ws.Range("G3").Value = 0

# Row 4
$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial(-4163)
$ws.Range("F4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 0
$ws.Range("G4").Value = 0

# Row 5
$ws.Range("G5").Copy()
$ws.Range("I5").PasteSpecial(-4163)
$ws.Range("F5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = 0
$ws.Range("G5").Value = 4956

# Row 6
$ws.Range("G6").Copy()
$ws.Range("I6").PasteSpecial(-4163)
$ws.Range("F6").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H6").Value = 0
$ws.Range("G6").Value = 2742

# Row 7
$ws.Range("G7").Copy()
$ws.Range("I7").PasteSpecial(-4163)
$ws.Range("F7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = 0
$ws.Range("G7").Value = 0

# Row 8
$ws.Range("G8").Copy()
$ws.Range("I8").PasteSpecial(-4163)
$ws.Range("F8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = 40
$ws.Range("G8").Value = 4975

# Row 9
$ws.Range("G9").Copy()
$ws.Range("I9").PasteSpecial(-4163)
$ws.Range("F9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = 5
$ws.Range("G9").Value = 4769

# Row 10
$ws.Range("G10").Copy()
$ws.Range("I10").PasteSpecial(-4163)
$ws.Range("F10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Value = 0
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("G11").Copy()
$ws.Range("I11").PasteSpecial(-4163)
$ws.Range("F11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = 0
$ws.Range("G11").Value = 2813

# Row 12
$ws.Range("G12").Copy()
$ws.Range("I12").PasteSpecial(-4163)
$ws.Range("F12").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = 0
$ws.Range("G12").Value = 1824

# Row 13
$ws.Range("G13").Copy()
$ws.Range("I13").PasteSpecial(-4163)
$ws.Range("F13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = 15
$ws.Range("G13").Value = 4198

# Row 14
$ws.Range("G14").Copy()
$ws.Range("I14").PasteSpecial(-4163)
$ws.Range("F14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = 0
$ws.Range("G14").Value = 0

# Row 15
$ws.Range("G15").Copy()
$ws.Range("I15").PasteSpecial(-4163)
$ws.Range("F15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 18
$ws.Range("G15").Value = 4621

# Row 16
$ws.Range("G16").Copy()
$ws.Range("I16").PasteSpecial(-4163)
$ws.Range("F16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = 30
$ws.Range("G16").Value = 6537

# Row 17
$ws.Range("G17").Copy()
$ws.Range("I17").PasteSpecial(-4163)
$ws.Range("F17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H17").Value = 37
$ws.Range("G17").Value = 5402

# Row 18
$ws.Range("G18").Copy()
$ws.Range("I18").PasteSpecial(-4163)
$ws.Range("F18").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H18").Value = 21
$ws.Range("G18").Value = 5602

# Row 19
$ws.Range("G19").Copy()
$ws.Range("I19").PasteSpecial(-4163)
$ws.Range("F19").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value = 22
$ws.Range("G19").Value = 6503

# Row 20
$ws.Range("G20").Copy()
$ws.Range("I20").PasteSpecial(-4163)
$ws.Range("F20").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = 33
$ws.Range("G20").Value = 5665

# Row 21
$ws.Range("G21").Copy()
$ws.Range("I21").PasteSpecial(-4163)
$ws.Range("F21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H21").Value = 27
$ws.Range("G21").Value = 6176

# Row 22
$ws.Range("G22").Copy()
$ws.Range("I22").PasteSpecial(-4163)
$ws.Range("F22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = 16
$ws.Range("G22").Value = 5494

# Row 23
$ws.Range("G23").Copy()
$ws.Range("I23").PasteSpecial(-4163)
$ws.Range("F23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").Value = 20
$ws.Range("G23").Value = 4041

# Row 24
$ws.Range("G24").Copy()
$ws.Range("I24").PasteSpecial(-4163)
$ws.Range("F24").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("H24").Value = 0
$ws.Range("G24").Value = 4237

# Row 25
$ws.Range("G25").Copy()
$ws.Range("I25").PasteSpecial(-4163)
$ws.Range("F25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value = 30
$ws.Range("G25").Value = 5041

# Row 26
$ws.Range("G26").Copy()
$ws.Range("I26").PasteSpecial(-4163)
$ws.Range("F26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Value = 30
$ws.Range("G26").Value = 5579

# Row 27
$ws.Range("G27").Copy()
$ws.Range("I27").PasteSpecial(-4163)
$ws.Range("F27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = 21
$ws.Range("G27").Value = 4934

# Row 28
$ws.Range("G28").Copy()
$ws.Range("I28").PasteSpecial(-4163)
$ws.Range("F28").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = 0
$ws.Range("G28").Value = 5732

# Row 29
$ws.Range("G29").Copy()
$ws.Range("I29").PasteSpecial(-4163)
$ws.Range("F29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 20
$ws.Range("G29").Value = 5363

# Row 30
$ws.Range("G30").Copy()
$ws.Range("I30").PasteSpecial(-4163)
$ws.Range("F30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = 16
$ws.Range("G30").Value = 4610

# Row 31
$ws.Range("G31").Copy()
$ws.Range("I31").PasteSpecial(-4163)
$ws.Range("F31").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = 30
$ws.Range("G31").Value = 5176

# Row 32
$ws.Range("G32").Copy()
$ws.Range("I32").PasteSpecial(-4163)
$ws.Range("F32").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H32").Value = 28
$ws.Range("G32").Value = 4803

# Row 33
$ws.Range("G33").Copy()
$ws.Range("I33").PasteSpecial(-4163)
$ws.Range("F33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = 21
$ws.Range("G33").Value = 6303

# Row 34
$ws.Range("G34").Copy()
$ws.Range("I34").PasteSpecial(-4163)
$ws.Range("F34").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("H34").Value = 29
$ws.Range("G34").Value = 6095

# Row 35
$ws.Range("G35").Copy()
$ws.Range("I35").PasteSpecial(-4163)
$ws.Range("F35").Copy()
$ws.Range("H35").PasteSpecial(-4122)
$ws.Range("H35").Value = 0
$ws.Range("G35").Value = 4689

# Row 36
$ws.Range("G36").Copy()
$ws.Range("I36").PasteSpecial(-4163)
$ws.Range("F36").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("H36").Value = 0
$ws.Range("G36").Value = 3532

# Row 37
$ws.Range("G37").Copy()
$ws.Range("I37").PasteSpecial(-4163)
$ws.Range("F37").Copy()
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("H37").Value = 30
$ws.Range("G37").Value = 5973

# Row 38
$ws.Range("G38").Copy()
$ws.Range("I38").PasteSpecial(-4163)
$ws.Range("F38").Copy()
$ws.Range("H38").PasteSpecial(-4122)
$ws.Range("H38").Value = 27
$ws.Range("G38").Value = 4244

# Row 39
$ws.Range("G39").Copy()
$ws.Range("I39").PasteSpecial(-4163)
$ws.Range("F39").Copy()
$ws.Range("H39").PasteSpecial(-4122)
$ws.Range("H39").Value = 15
$ws.Range("G39").Value = 3861

# Row 40
$ws.Range("G40").Copy()
$ws.Range("I40").PasteSpecial(-4163)
$ws.Range("F40").Copy()
$ws.Range("H40").PasteSpecial(-4122)
$ws.Range("H40").Value = 33
$ws.Range("G40").Value = 6057

# Row 41
$ws.Range("G41").Copy()
$ws.Range("I41").PasteSpecial(-4163)
$ws.Range("F41").Copy()
$ws.Range("H41").PasteSpecial(-4122)
$ws.Range("H41").Value = 29
$ws.Range("G41").Value = 5869

# Row 42
$ws.Range("G42").Copy()
$ws.Range("I42").PasteSpecial(-4163)
$ws.Range("F42").Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("H42").Value = 24
$ws.Range("G42").Value = 5264

# Row 43
$ws.Range("G43").Copy()
$ws.Range("I43").PasteSpecial(-4163)
$ws.Range("F43").Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("H43").Value = 20
$ws.Range("G43").Value = 5532

# Row 44
$ws.Range("G44").Copy()
$ws.Range("I44").PasteSpecial(-4163)
$ws.Range("F44").Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("H44").Value = 2
$ws.Range("G44").Value = 5095

# Row 45
$ws.Range("G45").Copy()
$ws.Range("I45").PasteSpecial(-4163)
$ws.Range("F45").Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("H45").Value = 30
$ws.Range("G45").Value = 5438

# Row 46
$ws.Range("G46").Copy()
$ws.Range("I46").PasteSpecial(-4163)
$ws.Range("F46").Copy()
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("H46").Value = 23
$ws.Range("G46").Value = 5220

# Row 47
$ws.Range("G47").Copy()
$ws.Range("I47").PasteSpecial(-4163)
$ws.Range("F47").Copy()
$ws.Range("H47").PasteSpecial(-4122)
$ws.Range("H47").Value = 20
$ws.Range("G47").Value = 4936

# Row 48
$ws.Range("G48").Copy()
$ws.Range("I48").PasteSpecial(-4163)
$ws.Range("F48").Copy()
$ws.Range("H48").PasteSpecial(-4122)
$ws.Range("H48").Value = 30
$ws.Range("G48").Value = 5828

# Row 49
$ws.Range("G49").Copy()
$ws.Range("I49").PasteSpecial(-4163)
$ws.Range("F49").Copy()
$ws.Range("H49").PasteSpecial(-4122)
$ws.Range("H49").Value = 15
$ws.Range("G49").Value = 4556

# Row 50
$ws.Range("G50").Copy()
$ws.Range("I50").PasteSpecial(-4163)
$ws.Range("F50").Copy()
$ws.Range("H50").PasteSpecial(-4122)
$ws.Range("H50").Value = 5
$ws.Range("G50").Value = 4139

# Row 51
$ws.Range("G51").Copy()
$ws.Range("I51").PasteSpecial(-4163)
$ws.Range("F51").Copy()
$ws.Range("H51").PasteSpecial(-4122)
$ws.Range("H51").Value = 9
$ws.Range("G51").Value = 4616

# Row 52
$ws.Range("G52").Copy()
$ws.Range("I52").PasteSpecial(-4163)
$ws.Range("F52").Copy()
$ws.Range("H52").PasteSpecial(-4122)
$ws.Range("H52").Value = 18
$ws.Range("G52").Value = 4500

# Row 53
$ws.Range("G53").Copy()
$ws.Range("I53").PasteSpecial(-4163)
$ws.Range("F53").Copy()
$ws.Range("H53").PasteSpecial(-4122)
$ws.Range("H53").Value = 20
$ws.Range("G53").Value = 4596

# Row 54
$ws.Range("G54").Copy()
$ws.Range("I54").PasteSpecial(-4163)
$ws.Range("F54").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("H54").Value = 20
$ws.Range("G54").Value = 4055

# Row 55
$ws.Range("G55").Copy()
$ws.Range("I55").PasteSpecial(-4163)
$ws.Range("F55").Copy()
$ws.Range("H55").PasteSpecial(-4122)
$ws.Range("H55").Value = 30
$ws.Range("G55").Value = 4218

# Row 56
$ws.Range("G56").Copy()
$ws.Range("I56").PasteSpecial(-4163)
$ws.Range("F56").Copy()
$ws.Range("H56").PasteSpecial(-4122)
$ws.Range("H56").Value = 10
$ws.Range("G56").Value = 4123

# Row 57
$ws.Range("G57").Copy()
$ws.Range("I57").PasteSpecial(-4163)
$ws.Range("F57").Copy()
$ws.Range("H57").PasteSpecial(-4122)
$ws.Range("H57").Value = 17
$ws.Range("G57").Value = 4506

# Row 58
$ws.Range("G58").Copy()
$ws.Range("I58").PasteSpecial(-4163)
$ws.Range("F58").Copy()
$ws.Range("H58").PasteSpecial(-4122)
$ws.Range("H58").Value = 0
$ws.Range("G58").Value = 0

# Row 59
$ws.Range("G59").Copy()
$ws.Range("I59").PasteSpecial(-4163)
$ws.Range("F59").Copy()
$ws.Range("H59").PasteSpecial(-4122)
$ws.Range("H59").Value = 0
$ws.Range("G59").Value = 2944

# Row 60
$ws.Range("G60").Copy()
$ws.Range("I60").PasteSpecial(-4163)
$ws.Range("F60").Copy()
$ws.Range("H60").PasteSpecial(-4122)
$ws.Range("H60").Value = 0
$ws.Range("G60").Value = 2737

# Row 61
$ws.Range("G61").Copy()
$ws.Range("I61").PasteSpecial(-4163)
$ws.Range("F61").Copy()
$ws.Range("H61").PasteSpecial(-4122)
$ws.Range("H61").Value = 0
$ws.Range("G61").Value = 3168

# Row 62
$ws.Range("G62").Copy()
$ws.Range("I62").PasteSpecial(-4163)
$ws.Range("F62").Copy()
$ws.Range("H62").PasteSpecial(-4122)
$ws.Range("H62").Value = 0
$ws.Range("G62").Value = 0

# Row 63
$ws.Range("F63").Copy()
$ws.Range("H63").PasteSpecial(-4122)
$ws.Range("G63").Copy()
$ws.Range("I63").PasteSpecial(-4122)

# Row 64
$ws.Range("G64").Copy()
$ws.Range("I64").PasteSpecial(-4163)
$ws.Range("F64").Copy()
$ws.Range("H64").PasteSpecial(-4122)
$ws.Range("H64").Value = 0
$ws.Range("G64").Value = 0

# Row 65
$ws.Range("G65").Copy()
$ws.Range("I65").PasteSpecial(-4163)
$ws.Range("F65").Copy()
$ws.Range("H65").PasteSpecial(-4122)
$ws.Range("H65").Value = 20
$ws.Range("G65").Value = 5200

# Row 66
$ws.Range("G66").Copy()
$ws.Range("I66").PasteSpecial(-4163)
$ws.Range("F66").Copy()
$ws.Range("H66").PasteSpecial(-4122)
$ws.Range("H66").Value = 0
$ws.Range("G66").Value = 0

# Row 67
$ws.Range("G67").Copy()
$ws.Range("I67").PasteSpecial(-4163)
$ws.Range("F67").Copy()
$ws.Range("H67").PasteSpecial(-4122)
$ws.Range("H67").Value = 28
$ws.Range("G67").Value = 5509

# Row 68
$ws.Range("G68").Copy()
$ws.Range("I68").PasteSpecial(-4163)
$ws.Range("F68").Copy()
$ws.Range("H68").PasteSpecial(-4122)
$ws.Range("H68").Value = 0
$ws.Range("G68").Value = 0

# Row 69
$ws.Range("G69").Copy()
$ws.Range("I69").PasteSpecial(-4163)
$ws.Range("F69").Copy()
$ws.Range("H69").PasteSpecial(-4122)
$ws.Range("H69").Value = 0
$ws.Range("G69").Value = 2757

# Row 70
$ws.Range("G70").Copy()
$ws.Range("I70").PasteSpecial(-4163)
$ws.Range("F70").Copy()
$ws.Range("H70").PasteSpecial(-4122)
$ws.Range("H70").Value = 0
$ws.Range("G70").Value = 0

# Row 71
$ws.Range("G71").Copy()
$ws.Range("I71").PasteSpecial(-4163)
$ws.Range("F71").Copy()
$ws.Range("H71").PasteSpecial(-4122)
$ws.Range("H71").Value = 2
$ws.Range("G71").Value = 4891

# Row 72
$ws.Range("G72").Copy()
$ws.Range("I72").PasteSpecial(-4163)
$ws.Range("F72").Copy()
$ws.Range("H72").PasteSpecial(-4122)
$ws.Range("H72").Value = 0
$ws.Range("G72").Value = 3355

# Row 73
$ws.Range("G73").Copy()
$ws.Range("I73").PasteSpecial(-4163)
$ws.Range("F73").Copy()
$ws.Range("H73").PasteSpecial(-4122)
$ws.Range("H73").Value = 0
$ws.Range("G73").Value = 0

# Row 74
$ws.Range("G74").Copy()
$ws.Range("I74").PasteSpecial(-4163)
$ws.Range("F74").Copy()
$ws.Range("H74").PasteSpecial(-4122)
$ws.Range("H74").Value = 0
$ws.Range("G74").Value = 1225

# Row 75
$ws.Range("G75").Copy()
$ws.Range("I75").PasteSpecial(-4163)
$ws.Range("F75").Copy()
$ws.Range("H75").PasteSpecial(-4122)
$ws.Range("H75").Value = 0
$ws.Range("G75").Value = 3440

# Row 76
$ws.Range("G76").Copy()
$ws.Range("I76").PasteSpecial(-4163)
$ws.Range("F76").Copy()
$ws.Range("H76").PasteSpecial(-4122)
$ws.Range("H76").Value = 0
$ws.Range("G76").Value = 0

# Row 77
$ws.Range("G77").Copy()
$ws.Range("I77").PasteSpecial(-4163)
$ws.Range("F77").Copy()
$ws.Range("H77").PasteSpecial(-4122)
$ws.Range("H77").Value = 20
$ws.Range("G77").Value = 4653

# Row 78
$ws.Range("G78").Copy()
$ws.Range("I78").PasteSpecial(-4163)
$ws.Range("F78").Copy()
$ws.Range("H78").PasteSpecial(-4122)
$ws.Range("H78").Value = 0
$ws.Range("G78").Value = 0

# Row 79
$ws.Range("G79").Copy()
$ws.Range("I79").PasteSpecial(-4163)
$ws.Range("F79").Copy()
$ws.Range("H79").PasteSpecial(-4122)
$ws.Range("H79").Value = 0
$ws.Range("G79").Value = 3187

# Row 80
$ws.Range("G80").Copy()
$ws.Range("I80").PasteSpecial(-4163)
$ws.Range("F80").Copy()
$ws.Range("H80").PasteSpecial(-4122)
$ws.Range("H80").Value = 0
$ws.Range("G80").Value = 1455

# Row 81
$ws.Range("G81").Copy()
$ws.Range("I81").PasteSpecial(-4163)
$ws.Range("F81").Copy()
$ws.Range("H81").PasteSpecial(-4122)
$ws.Range("H81").Value = 0
$ws.Range("G81").Value = 0

# Row 82
$ws.Range("G82").Copy()
$ws.Range("I82").PasteSpecial(-4163)
$ws.Range("F82").Copy()
$ws.Range("H82").PasteSpecial(-4122)
$ws.Range("H82").Value = 0
$ws.Range("G82").Value = 0

# Row 83
$ws.Range("G83").Copy()
$ws.Range("I83").PasteSpecial(-4163)
$ws.Range("F83").Copy()
$ws.Range("H83").PasteSpecial(-4122)
$ws.Range("H83").Value = 0
$ws.Range("G83").Value = 0

# Row 84
$ws.Range("G84").Copy()
$ws.Range("I84").PasteSpecial(-4163)
$ws.Range("F84").Copy()
$ws.Range("H84").PasteSpecial(-4122)
$ws.Range("H84").Value = 0
$ws.Range("G84").Value = 0

# Row 85
$ws.Range("G85").Copy()
$ws.Range("I85").PasteSpecial(-4163)
$ws.Range("F85").Copy()
$ws.Range("H85").PasteSpecial(-4122)
$ws.Range("H85").Value = 0
$ws.Range("G85").Value = 0

# Row 86
$ws.Range("G86").Copy()
$ws.Range("I86").PasteSpecial(-4163)
$ws.Range("F86").Copy()
$ws.Range("H86").PasteSpecial(-4122)
$ws.Range("H86").Value = 0
$ws.Range("G86").Value = 999

# Row 87
$ws.Range("G87").Copy()
$ws.Range("I87").PasteSpecial(-4163)
$ws.Range("F87").Copy()
$ws.Range("H87").PasteSpecial(-4122)
$ws.Range("H87").Value = 0
$ws.Range("G87").Value = 0

# Row 88
$ws.Range("G88").Copy()
$ws.Range("I88").PasteSpecial(-4163)
$ws.Range("F88").Copy()
$ws.Range("H88").PasteSpecial(-4122)
$ws.Range("H88").Value = 0
$ws.Range("G88").Value = 3122

# Row 89
$ws.Range("G89").Copy()
$ws.Range("I89").PasteSpecial(-4163)
$ws.Range("F89").Copy()
$ws.Range("H89").PasteSpecial(-4122)
$ws.Range("H89").Value = 0
$ws.Range("G89").Value = 0

# Row 90
$ws.Range("G90").Copy()
$ws.Range("I90").PasteSpecial(-4163)
$ws.Range("F90").Copy()
$ws.Range("H90").PasteSpecial(-4122)
$ws.Range("H90").Value = 0
$ws.Range("G90").Value = 0

# Row 91
$ws.Range("G91").Copy()
$ws.Range("I91").PasteSpecial(-4163)
$ws.Range("F91").Copy()
$ws.Range("H91").PasteSpecial(-4122)
$ws.Range("H91").Value = 0
$ws.Range("G91").Value = 2899

# Row 92
$ws.Range("G92").Copy()
$ws.Range("I92").PasteSpecial(-4163)
$ws.Range("F92").Copy()
$ws.Range("H92").PasteSpecial(-4122)
$ws.Range("H92").Value = 0
$ws.Range("G92").Value = 0

# Row 93
$ws.Range("G93").Copy()
$ws.Range("I93").PasteSpecial(-4163)
$ws.Range("F93").Copy()
$ws.Range("H93").PasteSpecial(-4122)
$ws.Range("H93").Value = 0
$ws.Range("G93").Value = 3982

# Row 94
$ws.Range("G94").Copy()
$ws.Range("I94").PasteSpecial(-4163)
$ws.Range("F94").Copy()
$ws.Range("H94").PasteSpecial(-4122)
$ws.Range("H94").Value = 0
$ws.Range("G94").Value = 0

# Row 95
$ws.Range("G95").Copy()
$ws.Range("I95").PasteSpecial(-4163)
$ws.Range("F95").Copy()
$ws.Range("H95").PasteSpecial(-4122)
$ws.Range("H95").Value = 0
$ws.Range("G95").Value = 0

# Row 96
$ws.Range("G96").Copy()
$ws.Range("I96").PasteSpecial(-4163)
$ws.Range("F96").Copy()
$ws.Range("H96").PasteSpecial(-4122)
$ws.Range("H96").Value = 0
$ws.Range("G96").Value = 2057

# Row 97
$ws.Range("G97").Copy()
$ws.Range("I97").PasteSpecial(-4163)
$ws.Range("F97").Copy()
$ws.Range("H97").PasteSpecial(-4122)
$ws.Range("H97").Value = 0
$ws.Range("G97").Value = 0

# Row 98
$ws.Range("G98").Copy()
$ws.Range("I98").PasteSpecial(-4163)
$ws.Range("F98").Copy()
$ws.Range("H98").PasteSpecial(-4122)
$ws.Range("H98").Value = 0
$ws.Range("G98").Value = 0

# Row 99
$ws.Range("G99").Copy()
$ws.Range("I99").PasteSpecial(-4163)
$ws.Range("F99").Copy()
$ws.Range("H99").PasteSpecial(-4122)
$ws.Range("H99").Value = 0
$ws.Range("G99").Value = 0

# Row 100
$ws.Range("G100").Copy()
$ws.Range("I100").PasteSpecial(-4163)
$ws.Range("F100").Copy()
$ws.Range("H100").PasteSpecial(-4122)
$ws.Range("H100").Value = 0
$ws.Range("G100").Value = 2059

# Row 101
$ws.Range("G101").Copy()
$ws.Range("I101").PasteSpecial(-4163)
$ws.Range("F101").Copy()
$ws.Range("H101").PasteSpecial(-4122)
$ws.Range("H101").Value = 0
$ws.Range("G101").Value = 0

# Row 102
$ws.Range("G102").Copy()
$ws.Range("I102").PasteSpecial(-4163)
$ws.Range("F102").Copy()
$ws.Range("H102").PasteSpecial(-4122)
$ws.Range("H102").Value = 0
$ws.Range("G102").Value = 0

# Row 103
$ws.Range("G103").Copy()
$ws.Range("I103").PasteSpecial(-4163)
$ws.Range("F103").Copy()
$ws.Range("H103").PasteSpecial(-4122)
$ws.Range("H103").Value = 0
$ws.Range("G103").Value = 0

# Row 104
$ws.Range("G104").Copy()
$ws.Range("I104").PasteSpecial(-4163)
$ws.Range("F104").Copy()
$ws.Range("H104").PasteSpecial(-4122)
$ws.Range("H104").Value = 0
$ws.Range("G104").Value = 0

# Row 105
$ws.Range("G105").Copy()
$ws.Range("I105").PasteSpecial(-4163)
$ws.Range("F105").Copy()
$ws.Range("H105").PasteSpecial(-4122)
$ws.Range("H105").Value = 0
$ws.Range("G105").Value = 0

# Row 106
$ws.Range("G106").Copy()
$ws.Range("I106").PasteSpecial(-4163)
$ws.Range("F106").Copy()
$ws.Range("H106").PasteSpecial(-4122)
$ws.Range("H106").Value = 0
$ws.Range("G106").Value = 0

# Row 107
$ws.Range("G107").Copy()
$ws.Range("I107").PasteSpecial(-4163)
$ws.Range("F107").Copy()
$ws.Range("H107").PasteSpecial(-4122)
$ws.Range("H107").Value = 0
$ws.Range("G107").Value = 0

# Row 108
$ws.Range("G108").Copy()
$ws.Range("I108").PasteSpecial(-4163)
$ws.Range("F108").Copy()
$ws.Range("H108").PasteSpecial(-4122)
$ws.Range("H108").Value = 0
$ws.Range("G108").Value = 0

# Row 109
$ws.Range("G109").Copy()
$ws.Range("I109").PasteSpecial(-4163)
$ws.Range("F109").Copy()
$ws.Range("H109").PasteSpecial(-4122)
$ws.Range("H109").Value = 0
$ws.Range("G109").Value = 0

# Row 110
$ws.Range("G110").Copy()
$ws.Range("I110").PasteSpecial(-4163)
$ws.Range("F110").Copy()
$ws.Range("H110").PasteSpecial(-4122)
$ws.Range("H110").Value = 0
$ws.Range("G110").Value = 0

# Row 111
$ws.Range("G111").Copy()
$ws.Range("I111").PasteSpecial(-4163)
$ws.Range("F111").Copy()
$ws.Range("H111").PasteSpecial(-4122)
$ws.Range("H111").Value = 0
$ws.Range("G111").Value = 3224
